$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell address -> new text value. Column D entries that look like
# plain numbers are flagged so we can force them to stay text cells (they
# are price strings such as "8.76", "43.791.43", not numeric cells).
$updates = @(
    @{ Cell = 'D2'; Value = '43.791.43'; ForceText = 0 }
    @{ Cell = 'E2'; Value = '  +0.47%  '; ForceText = 0 }
    @{ Cell = 'D3'; Value = '2.306.22'; ForceText = 0 }
    @{ Cell = 'E3'; Value = '  +1.00%  '; ForceText = 0 }
    @{ Cell = 'E4'; Value = '  +0.14%  '; ForceText = 0 }
    @{ Cell = 'D5'; Value = '114.68'; ForceText = 1 }
    @{ Cell = 'E5'; Value = '  +20.91%  '; ForceText = 0 }
    @{ Cell = 'D6'; Value = '269.44'; ForceText = 1 }
    @{ Cell = 'E6'; Value = '  +1.10%  '; ForceText = 0 }
    @{ Cell = 'E7'; Value = '  +1.05%  '; ForceText = 0 }
    @{ Cell = 'E8'; Value = '  +0.24%  '; ForceText = 0 }
    @{ Cell = 'D9'; Value = '0.626'; ForceText = 1 }
    @{ Cell = 'E9'; Value = '  +3.39%  '; ForceText = 0 }
    @{ Cell = 'D10'; Value = '48.56'; ForceText = 1 }
    @{ Cell = 'E10'; Value = '  +9.13%  '; ForceText = 0 }
    @{ Cell = 'D11'; Value = '0.0945'; ForceText = 1 }
    @{ Cell = 'E11'; Value = '  +1.83%  '; ForceText = 0 }
    @{ Cell = 'D12'; Value = '8.76'; ForceText = 1 }
    @{ Cell = 'E12'; Value = '  +13.36%  '; ForceText = 0 }
    @{ Cell = 'E13'; Value = '  +2.87%  '; ForceText = 0 }
    @{ Cell = 'D14'; Value = '15.64'; ForceText = 1 }
    @{ Cell = 'E14'; Value = '  +3.48%  '; ForceText = 0 }
    @{ Cell = 'D15'; Value = '2.620.14'; ForceText = 0 }
    @{ Cell = 'E15'; Value = '  -0.17%  '; ForceText = 0 }
    @{ Cell = 'E16'; Value = '  +2.14%  '; ForceText = 0 }
    @{ Cell = 'D17'; Value = '2.307.09'; ForceText = 0 }
    @{ Cell = 'E17'; Value = '  +1.01%  '; ForceText = 0 }
    @{ Cell = 'D18'; Value = '43.762.68'; ForceText = 0 }
    @{ Cell = 'E18'; Value = '  +0.57%  '; ForceText = 0 }
    @{ Cell = 'D19'; Value = '0.0000111'; ForceText = 1 }
    @{ Cell = 'E19'; Value = '  +4.17%  '; ForceText = 0 }
    @{ Cell = 'D20'; Value = '6.60'; ForceText = 1 }
    @{ Cell = 'E20'; Value = '  +7.17%  '; ForceText = 0 }
    @{ Cell = 'D21'; Value = '72.66'; ForceText = 1 }
    @{ Cell = 'E21'; Value = '  +0.68%  '; ForceText = 0 }
    @{ Cell = 'D22'; Value = '2.57'; ForceText = 1 }
    @{ Cell = 'E22'; Value = '  +7.70%  '; ForceText = 0 }
    @{ Cell = 'D23'; Value = '234.73'; ForceText = 1 }
    @{ Cell = 'E23'; Value = '  +0.34%  '; ForceText = 0 }
    @{ Cell = 'D24'; Value = '9.55'; ForceText = 1 }
    @{ Cell = 'E24'; Value = '  +6.97%  '; ForceText = 0 }
    @{ Cell = 'D25'; Value = '2.86'; ForceText = 1 }
    @{ Cell = 'E25'; Value = '  +14.96%  '; ForceText = 0 }
    @{ Cell = 'E26'; Value = '  -0.01%  '; ForceText = 0 }
    @{ Cell = 'D27'; Value = '11.56'; ForceText = 1 }
    @{ Cell = 'E27'; Value = '  +4.05%  '; ForceText = 0 }
    @{ Cell = 'D28'; Value = '43.57'; ForceText = 1 }
    @{ Cell = 'E28'; Value = '  +10.43%  '; ForceText = 0 }
    @{ Cell = 'E29'; Value = '  -1.31%  '; ForceText = 0 }
    @{ Cell = 'D31'; Value = '178.08'; ForceText = 1 }
    @{ Cell = 'D32'; Value = '21.93'; ForceText = 1 }
    @{ Cell = 'E32'; Value = '  +0.46%  '; ForceText = 0 }
    @{ Cell = 'D33'; Value = '0.0935'; ForceText = 1 }
    @{ Cell = 'E33'; Value = '  +6.50%  '; ForceText = 0 }
    @{ Cell = 'E34'; Value = '  +4.57%  '; ForceText = 0 }
    @{ Cell = 'E35'; Value = '  +1.40%  '; ForceText = 0 }
    @{ Cell = 'E36'; Value = '  +8.41%  '; ForceText = 0 }
    @{ Cell = 'E37'; Value = '  +3.23%  '; ForceText = 0 }
    @{ Cell = 'D38'; Value = '3.96'; ForceText = 1 }
    @{ Cell = 'E38'; Value = '  +20.06%  '; ForceText = 0 }
    @{ Cell = 'D39'; Value = '0.0357'; ForceText = 1 }
    @{ Cell = 'E39'; Value = '  +1.09%  '; ForceText = 0 }
    @{ Cell = 'D40'; Value = '75.64'; ForceText = 1 }
    @{ Cell = 'E40'; Value = '  +19.12%  '; ForceText = 0 }
    @{ Cell = 'E41'; Value = '  +5.55%  '; ForceText = 0 }
    @{ Cell = 'D42'; Value = '2.41'; ForceText = 1 }
    @{ Cell = 'E42'; Value = '  +3.14%  '; ForceText = 0 }
    @{ Cell = 'D43'; Value = '13.38'; ForceText = 1 }
    @{ Cell = 'E43'; Value = '  +12.53%  '; ForceText = 0 }
    @{ Cell = 'B44'; Value = 'ARBITRUM'; ForceText = 0 }
    @{ Cell = 'C44'; Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'; ForceText = 0 }
    @{ Cell = 'D44'; Value = '1.42'; ForceText = 1 }
    @{ Cell = 'E44'; Value = '  +6.29%  '; ForceText = 0 }
    @{ Cell = 'B45'; Value = 'THORChain'; ForceText = 0 }
    @{ Cell = 'C45'; Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'; ForceText = 0 }
    @{ Cell = 'D45'; Value = '6.01'; ForceText = 1 }
    @{ Cell = 'E45'; Value = '  +15.35%  '; ForceText = 0 }
    @{ Cell = 'B46'; Value = 'FirstDigitalUSD'; ForceText = 0 }
    @{ Cell = 'C46'; Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'; ForceText = 0 }
    @{ Cell = 'D46'; Value = '1.00'; ForceText = 1 }
    @{ Cell = 'E46'; Value = '  +0.08%  '; ForceText = 0 }
    @{ Cell = 'D47'; Value = '8.78'; ForceText = 1 }
    @{ Cell = 'E47'; Value = '  +0.11%  '; ForceText = 0 }
    @{ Cell = 'D49'; Value = '101.88'; ForceText = 1 }
    @{ Cell = 'E49'; Value = '  +4.48%  '; ForceText = 0 }
    @{ Cell = 'E50'; Value = '  +4.59%  '; ForceText = 0 }
    @{ Cell = 'D51'; Value = '0.460'; ForceText = 1 }
    @{ Cell = 'E51'; Value = '  +8.74%  '; ForceText = 0 }
)

foreach ($u in $updates) {
    $rng = $ws.Range($u.Cell)
    if ($u.ForceText -eq 1) {
        # Without this, Excel auto-converts numeric-looking text (e.g. "8.76")
        # into a real number when assigned via .Value. Force text storage, then
        # restore the default "Normal" style so no stray style index is left on
        # the cell.
        $rng.NumberFormat = "@"
        $rng.Value = $u.Value
        $rng.Style = "Normal"
    } else {
        $rng.Value = $u.Value
    }
}
